$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 2: "Defines procedures for using TWAMP Light messages for delay, synthetic loss
# and direct-mode loss measurements" -> "Defines procedures for delay, synthetic loss and
# direct-mode loss measurements"
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "."
$para2.Text = "Defines procedures for delay, synthetic loss and direct-mode loss measurements"

# Paragraph 3: demote to level 2 (OOXML lvl="1") and change text from
# "Procedures are defined for Links and end-to-end SR Paths for SR-MPLS and SRv6 data-planes"
# to "For Links and end-to-end SR Paths for SR-MPLS and SRv6 data planes"
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "."
$para3.Text = "For Links and end-to-end SR Paths for SR-MPLS and SRv6 data planes"
$para3.IndentLevel = 2

# Paragraph 7: demote to level 2 (OOXML lvl="1") and fix capitalization of "Code"
# "Defines Session-Sender Control code field for in-band response request" ->
# "Defines Session-Sender Control Code field for in-band response request"
$para7 = $tr.Paragraphs(7, 1)
$para7.Text = "."
$para7.Text = "Defines Session-Sender Control Code field for in-band response request"
$para7.IndentLevel = 2

# Paragraph 8: demote to level 2 (OOXML lvl="1"), text unchanged
$para8 = $tr.Paragraphs(8, 1)
$para8.IndentLevel = 2
